$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $ws.Range($cell).Value = "'" + $value
    $ws.Range($cell).Style = "Normal"
}

# Row 2
Set-TextCell 'D2' '29.927.03'
Set-TextCell 'E2' '  -0.05%  '

# Row 3
Set-TextCell 'D3' '1.906.38'
Set-TextCell 'E3' '  +0.37%  '

# Row 4
Set-TextCell 'D4' '1.002'
Set-TextCell 'E4' '  +0.30%  '

# Row 5
Set-TextCell 'D5' '0.7998'
Set-TextCell 'E5' '  +5.34%  '

# Row 6
Set-TextCell 'D6' '241.36'
Set-TextCell 'E6' '  +0.30%  '

# Row 7
Set-TextCell 'D7' '1.002'
Set-TextCell 'E7' '  +0.21%  '

# Row 8
Set-TextCell 'D8' '0.3138'
Set-TextCell 'E8' '  +2.59%  '

# Row 9
Set-TextCell 'D9' '26.23'
Set-TextCell 'E9' '  +2.74%  '

# Row 10
Set-TextCell 'D10' '0.06889'
Set-TextCell 'E10' '  +0.52%  '

# Row 11
Set-TextCell 'D11' '0.07987'
Set-TextCell 'E11' '  -0.02%  '

# Row 12
Set-TextCell 'D12' '1.915.75'
Set-TextCell 'E12' '  +0.90%  '

# Row 13
Set-TextCell 'D13' '0.7351'
Set-TextCell 'E13' '  -2.31%  '

# Row 14
Set-TextCell 'D14' '5.177'
Set-TextCell 'E14' '  -1.01%  '

# Row 15
Set-TextCell 'D15' '92.85'
Set-TextCell 'E15' '  +1.51%  '

# Row 16
Set-TextCell 'D16' '29.946.18'
Set-TextCell 'E16' '  +0.03%  '

# Row 17
Set-TextCell 'D17' '13.93'
Set-TextCell 'E17' '  -0.33%  '

# Row 18
Set-TextCell 'D18' '5.857'
Set-TextCell 'E18' '  -2.35%  '

# Row 19
Set-TextCell 'D19' '245.10'
Set-TextCell 'E19' '  +0.75%  '

# Row 20
Set-TextCell 'D20' '0.000007706'
Set-TextCell 'E20' '  -0.24%  '

# Row 21
Set-TextCell 'D21' '1.001'
Set-TextCell 'E21' '  +0.21%  '

# Row 22
Set-TextCell 'D22' '2.155.94'
Set-TextCell 'E22' '  +0.64%  '

# Row 23
Set-TextCell 'E23' '  +0.45%  '

# Row 24
Set-TextCell 'D24' '6.887'
Set-TextCell 'E24' '  -1.55%  '

# Row 25
Set-TextCell 'D25' '167.81'
Set-TextCell 'E25' '  +1.27%  '

# Row 26
Set-TextCell 'D26' '9.191'
Set-TextCell 'E26' '  -0.84%  '

# Row 27
Set-TextCell 'D27' '0.1416'
Set-TextCell 'E27' '  +8.99%  '

# Row 28
Set-TextCell 'D28' '18.86'
Set-TextCell 'E28' '  +0.38%  '

# Row 29
Set-TextCell 'D29' '2.023'
Set-TextCell 'E29' '  -0.96%  '

# Row 30
Set-TextCell 'E30' '  +0.45%  '

# Row 31
Set-TextCell 'D31' '1.513'
Set-TextCell 'E31' '  -0.35%  '

# Row 32
Set-TextCell 'D32' '4.292'
Set-TextCell 'E32' '  -0.27%  '

# Row 33
Set-TextCell 'D33' '4.063'
Set-TextCell 'E33' '  +0.59%  '

# Row 34
Set-TextCell 'D34' '0.05510'
Set-TextCell 'E34' '  +2.38%  '

# Row 35
Set-TextCell 'D35' '1.260'
Set-TextCell 'E35' '  +0.31%  '

# Row 36
Set-TextCell 'D36' '0.7295'
Set-TextCell 'E36' '  -0.20%  '

# Row 37
Set-TextCell 'D37' '2.723'
Set-TextCell 'E37' '  +0.05%  '

# Row 38
Set-TextCell 'D38' '0.01925'
Set-TextCell 'E38' '  -0.23%  '

# Row 39
Set-TextCell 'D39' '2.791'
Set-TextCell 'E39' '  +0.68%  '

# Row 40
Set-TextCell 'D40' '6.145'
Set-TextCell 'E40' '  -1.00%  '

# Row 41
Set-TextCell 'D41' '0.4403'
Set-TextCell 'E41' '  -0.50%  '

# Row 42
Set-TextCell 'D42' '72.08'
Set-TextCell 'E42' '  -0.64%  '

# Row 43
Set-TextCell 'E43' '  +0.23%  '

# Row 44
Set-TextCell 'D44' '0.8353'
Set-TextCell 'E44' '  +0.57%  '

# Row 45
Set-TextCell 'D45' '1.869'
Set-TextCell 'E45' '  -2.71%  '

# Row 46
Set-TextCell 'D46' '100.52'
Set-TextCell 'E46' '  -0.59%  '

# Row 47
Set-TextCell 'D47' '7.546'
Set-TextCell 'E47' '  -0.99%  '

# Row 48
Set-TextCell 'B48' 'EnergySwap'
Set-TextCell 'C48' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextCell 'D48' '9.720'
Set-TextCell 'E48' '  -0.47%  '

# Row 49
Set-TextCell 'D49' '2.063.58'
Set-TextCell 'E49' '  +0.68%  '

# Row 50
Set-TextCell 'B50' 'Maker'
Set-TextCell 'C50' 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell 'D50' '974.15'
Set-TextCell 'E50' '  +5.45%  '

# Row 51
Set-TextCell 'B51' 'Elrond'
Set-TextCell 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextCell 'D51' '36.16'
Set-TextCell 'E51' '  -0.35%  '
